# Update TPM-derived values in the LR-pairs sheet (Gdf6-Bmpr1a)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 4.993165333333334
$ws.Range("N2").Value = 14.979496
$ws.Range("O2").Value = 0.06779298131037136
$ws.Range("P2").Value = 0.06779298131037137
$ws.Range("Q2").Value = 1.798409995156445
$ws.Range("R2").Value = 16.185689956408
$ws.Range("S2").Value = 0.06779298131037136
$ws.Range("T2").Value = 0.06779298131037137

# Row 3
$ws.Range("O3").Value = 0.5355771637189464
$ws.Range("P3").Value = 0.5355771637189464
$ws.Range("S3").Value = 0.5355771637189464
$ws.Range("T3").Value = 0.5355771637189464

# Row 4
$ws.Range("M4").Value = 29.08216166666666
$ws.Range("N4").Value = 87.24648499999999
$ws.Range("O4").Value = 0.3948530262300277
$ws.Range("P4").Value = 0.3948530262300277
$ws.Range("Q4").Value = 10.47464819018389
$ws.Range("R4").Value = 94.27183371165498
$ws.Range("S4").Value = 0.3948530262300277
$ws.Range("T4").Value = 0.3948530262300277

# Row 5
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.130869
$ws.Range("N5").Value = 0.392607
$ws.Range("O5").Value = 0.001776828740654623
$ws.Range("P5").Value = 0.001776828740654624
$ws.Range("Q5").Value = 0.047135654829
$ws.Range("R5").Value = 0.424220893461
$ws.Range("S5").Value = 0.001776828740654623
$ws.Range("T5").Value = 0.001776828740654624
